# Fix markdown rendering: the category-label cells in column E stored their
# line separators as literal "<br>" tag text (which doesn't render as a line
# break since the cells aren't HTML). Replace each "<br>" with an actual
# line-break character so the text wraps onto multiple lines in Excel.
#
# Every cell sharing one of the three labels must be updated, since they all
# point at the same shared-string entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

$global = "**Global**: ${nl}Implemented by ${nl}All other countries"
$highIncome = "**High-income**: ${nl}All other HICs and ${nl}not some MICs (such as China)"
$international = "**International**: ${nl}Some countries (e.g. EU, UK, Brazil) ${nl}and not others (e.g. U.S., China)"

for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 5).Value = $global
}
for ($r = 14; $r -le 25; $r++) {
    $ws.Cells.Item($r, 5).Value = $highIncome
}
for ($r = 26; $r -le 37; $r++) {
    $ws.Cells.Item($r, 5).Value = $international
}
